# Fix database schema migration - append rows 144-146 to each of the four
# worksheets (MID_LFT_#1, MID_LFT_#2, MID_PLT_#1, MID_PLT_#2), continuing
# the daily log in the same style as the existing rows.

$wb = $excel.ActiveWorkbook

function Add-DataRow {
    param($ws, $row, $a, $b, $c, $d, $e, $f, $g, $h, $i)

    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = [double]$g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

# --- Sheet 1: MID_LFT_#1 ---
$ws1 = $wb.Worksheets.Item(1)
Add-DataRow $ws1 144 45930.46016203704 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x00,0xD4" "0x07" 400 "5.68631262647113e+23" 220 7
Add-DataRow $ws1 145 45931.4603125     "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x00,0xD4" "0x07" 400 "5.68631262647113e+23" 220 7
Add-DataRow $ws1 146 45932.46620370371 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x00,0xD4" "0x07" 400 "5.68631262647113e+23" 216 7

# --- Sheet 2: MID_LFT_#2 ---
$ws2 = $wb.Worksheets.Item(2)
Add-DataRow $ws2 144 45930.46016203704 "0x01,0x7c" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x00,0xF0" "0x19" 380 "5.68432987514711e+23" 240 25
Add-DataRow $ws2 145 45931.4603125     "0x01,0x7c" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x00,0xF0" "0x19" 380 "5.68432987514711e+23" 240 25
Add-DataRow $ws2 146 45932.46620370371 "0x01,0x7c" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x00,0xF0" "0x19" 380 "5.68432987514711e+23" 240 25

# --- Sheet 3: MID_PLT_#1 ---
$ws3 = $wb.Worksheets.Item(3)
Add-DataRow $ws3 144 45930.46016203704 "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x00,0x54" "0x15" 110 "5.68631262647113e+23" 84 15
Add-DataRow $ws3 145 45931.4603125     "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x00,0x54" "0x15" 110 "5.68631262647113e+23" 84 15
Add-DataRow $ws3 146 45932.46620370371 "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x00,0x53" "0x15" 110 "5.68631262647113e+23" 83 15

# --- Sheet 4: MID_PLT_#2 ---
$ws4 = $wb.Worksheets.Item(4)
Add-DataRow $ws4 144 45930.46016203704 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x00,0x6C" "0x9" 130 "5.68631262647113e+23" 108 9
Add-DataRow $ws4 145 45931.4603125     "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x00,0x6C" "0x9" 130 "5.68631262647113e+23" 108 9
Add-DataRow $ws4 146 45932.46620370371 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x00,0x6B" "0x9" 130 "5.68631262647113e+23" 107 9
